$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel stores them as literal text (matching the source sheet's inline strings)
# instead of silently coercing them to numbers.

$ws.Range("D2").Value = "60.693.27"
$ws.Range("E2").Value = "  +2.94%  "

$ws.Range("D3").Value = "2.602.61"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'569.37"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").Value = "'142.76"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").Value = "2.625.18"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("D10").Value = "'6.49"
$ws.Range("E10").Value = "  -2.74%  "

$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").Value = "'0.154"
$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("D13").Value = "'0.369"
$ws.Range("E13").Value = "  +7.29%  "

$ws.Range("D14").Value = "3.068.22"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "60.715.10"
$ws.Range("E15").Value = "  +2.92%  "

$ws.Range("D16").Value = "'23.59"
$ws.Range("E16").Value = "  +5.38%  "

$ws.Range("E17").Value = "  +3.10%  "

$ws.Range("D18").Value = "2.614.35"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19: coin identity change
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'11.29"
$ws.Range("E19").Value = "  +10.67%  "

# Row 20: coin identity change
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.68"
$ws.Range("E20").Value = "  +3.46%  "

$ws.Range("D21").Value = "'346.65"
$ws.Range("E21").Value = "  +3.58%  "

$ws.Range("D22").Value = "'7.12"
$ws.Range("E22").Value = "  +14.80%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'0.520"
$ws.Range("E24").Value = "  +14.20%  "

$ws.Range("D25").Value = "'63.69"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("E28").Value = "  +6.66%  "

$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").Value = "'1.80"
$ws.Range("E30").Value = "  +7.47%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").Value = "'161.50"
$ws.Range("E33").Value = "  +1.94%  "

$ws.Range("D34").Value = "'19.48"
$ws.Range("E34").Value = "  +2.78%  "

$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  +5.89%  "

$ws.Range("D36").Value = "'0.955"
$ws.Range("E36").Value = "  +9.59%  "

$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  +6.03%  "

$ws.Range("D39").Value = "'37.72"
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("D41").Value = "'3.79"
$ws.Range("E41").Value = "  +3.52%  "

$ws.Range("D42").Value = "'296.52"
$ws.Range("E42").Value = "  +1.43%  "

$ws.Range("D43").Value = "'138.72"
$ws.Range("E43").Value = "  +10.92%  "

$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "'0.0984"
$ws.Range("E45").Value = "  +0.90%  "

$ws.Range("E46").Value = "  +3.18%  "

$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("E48").Value = "  +3.59%  "

$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").Value = "'19.71"
$ws.Range("E50").Value = "  +7.07%  "

# Row 51: coin identity change
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.81"
$ws.Range("E51").Value = "  +7.12%  "
